$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new sprite reference strings for column B (rows 2-9), linking the
# item database to the inventory system's sprite names. Order matters here
# so the shared-strings table is built in the same sequence as the target.
$ws.Range("B2").Value = "spr_health_potion_light"
$ws.Range("B6").Value = "spr_armor_potion_light"
$ws.Range("B3").Value = "spr_placeholder"
$ws.Range("B4").Value = "spr_placeholder"
$ws.Range("B5").Value = "spr_placeholder"
$ws.Range("B7").Value = "spr_placeholder"
$ws.Range("B8").Value = "spr_placeholder"
$ws.Range("B9").Value = "spr_placeholder"

# Update the active cell selection to match the new edit location.
$ws.Range("C14").Select()
